$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Make room for the two new columns that will hold the new
#    'DOUBLE' cell-type fixtures (string "1.1" / "1,1" pair + the
#    formula that rebuilds each of them). Inserting at M:N pushes the
#    old M:Q range (date / currency / String / String / Automatic)
#    two columns to the right, landing on O:S.
# ------------------------------------------------------------------
$null = $ws.Range("M1:N1").EntireColumn.Insert()

# ------------------------------------------------------------------
# 2. G1 keeps its existing "@" text format and value of 1 - nothing
#    to change there. H1 used to read G1's text twice via a formula;
#    it becomes a literal CONCAT(1,1) with default ("Normal") styling
#    (no explicit number format anymore).
# ------------------------------------------------------------------
$ws.Range("H1").Style = "Normal"
$ws.Range("H1").Formula = "=CONCAT(1,1)"

# ------------------------------------------------------------------
# 3. K1/L1 used to hold the numeric 1.1 (stored as text) plus a
#    formula proving CONCAT(1,",1") renders as "1,1". They now hold
#    the "1.1" textual double plus the formula that rebuilds "1.1".
# ------------------------------------------------------------------
$ws.Range("K1").NumberFormat = "@"
$ws.Range("K1").Value = "1.1"
$ws.Range("L1").Style = "Normal"
$ws.Range("L1").Formula = '=CONCAT("1.",1)'

# ------------------------------------------------------------------
# 4. M1/N1 are the brand-new pair: the "1,1" textual double plus the
#    formula that rebuilds "1,1".
# ------------------------------------------------------------------
$ws.Range("M1").NumberFormat = "@"
$ws.Range("M1").Value = "1,1"
$ws.Range("N1").Style = "Normal"
$ws.Range("N1").Formula = '=CONCAT("1,",1)'

# ------------------------------------------------------------------
# 5. The cell comment ("Note") used to sit on the String cell that
#    was O1; that cell is now Q1. Column inserts don't relocate
#    comments automatically, so move it by hand.
# ------------------------------------------------------------------
$oldComment = $ws.Range("O1").Comment
$commentText = $oldComment.Text()
$oldComment.Delete()
$null = $ws.Range("Q1").AddComment($commentText)

# ------------------------------------------------------------------
# 6. Match the new active selection.
# ------------------------------------------------------------------
$null = $ws.Range("P1").Select()
